$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = '''source_file'
$ws.Range("B1").Value = '''Latitude'
$ws.Range("C1").Value = '''Longitude'
$ws.Range("D1").Value = '''Coordinates'
$ws.Range("E1").Value = '''Date_of_Issuance'
$ws.Range("F1").Value = '''Condition_summary_1'
$ws.Range("G1").Value = '''Condition_summary_2'
$ws.Range("H1").Value = '''Condition_summary_3'
$ws.Range("I1").Value = '''Condition_summary_4'
$ws.Range("J1").Value = '''Habitat_Type'
$ws.Range("K1").Value = '''Fish_species'
$ws.Range("L1").Value = '''Offset_footprint_size'
$ws.Range("M1").Value = '''Vegetation_Cover'
$ws.Range("N1").Value = '''Boulder'
$ws.Range("O1").Value = '''Woody_coverage'
$ws.Range("P1").Value = '''Instream_structures'

$ws.Range("A2").Value = '''OCR_18-HCAA-00233.json'
$ws.Range("B2").Value = '''None'
$ws.Range("C2").Value = '''None'
$ws.Range("D2").Value = '''None'
$ws.Range("E2").Value = '''JAN 10 2020'
$ws.Range("F2").Value = '''Install vegetation-enhanced armour stone walls on a 25 m section of the east bank (50 m²) and retrofit 10 parking lot catchbasins with shields along a 565 m section of the Moira River east bank.'
$ws.Range("G2").Value = '''Monitor offsetting measures for two years post-construction, reporting annually by November 30.'
$ws.Range("H2").Value = '''Implement contingency measures if offsetting measures fail to meet criteria (e.g., 80% vegetation survival or 50% turbidity reduction).'
$ws.Range("I2").Value = '''Prohibit adverse impacts on offsetting measures and ensure DFO access rights for implementation and monitoring.'
$ws.Range("J2").Value = '''Riverine'
$ws.Range("K2").Value = '''Channel Darter'
$ws.Range("L2").Value = '''50 m² vegetation-enhanced stone revetment, 10 catchbasin shields over 565 m river section'
$ws.Range("M2").Value = '''Riparian vegetation (enhanced through stone revetment)'
$ws.Range("N2").Value = '''None explicitly mentioned'
$ws.Range("O2").Value = '''None explicitly mentioned'
$ws.Range("P2").Value = '''Vegetation-enhanced armour stone walls (structural habitat), catchbasin shields (turbidity reduction infrastructure)'

$ws.Range("A3").Value = '''OCR_14-HCAA-00814.json'
$ws.Range("B3").Value = '''None'
$ws.Range("C3").Value = '''None'
$ws.Range("D3").Value = '''None'
$ws.Range("E3").Value = '''JUN 05 2015'
$ws.Range("F3").Value = '''Annual dam operation monitoring reports must be submitted by August 31 each year, summarizing dam operations from October 1 to July 15, including deviations from the operating plan and recommendations for changes.'
$ws.Range("G3").Value = '''Annual reports on juvenile Lake Sturgeon and Lake Whitefish recruitment must be submitted, detailing capture data, locations, substrates, depths, and age assessments using specified methods.'
$ws.Range("H3").Value = '''Contingency monitoring reports for additional offsetting measures must be provided as directed by DFO, describing adaptive management and effectiveness.'
$ws.Range("I3").Value = '''A final summary report combining all monitoring reports must be submitted within 90 days after the last report.'
$ws.Range("J3").Value = '''Riverine'
$ws.Range("K3").Value = '''Lake Sturgeon, Lake Whitefish'
$ws.Range("L3").Value = '''None'
$ws.Range("M3").Value = '''None'
$ws.Range("N3").Value = '''None'
$ws.Range("O3").Value = '''None'
$ws.Range("P3").Value = '''None'

$ws.Range("A4").Value = '''OCR_14-HCAA-00810.json'
$ws.Range("B4").Value = '''None'
$ws.Range("C4").Value = '''None'
$ws.Range("D4").Value = '''None'
$ws.Range("E4").Value = '''may 08 2015'
$ws.Range("F4").Value = '''offsetting measures must be completed during construction phase and meet criteria in the proponent’s plan'
$ws.Range("G4").Value = '''annual monitoring reports required, including fish rescue details and photos'
$ws.Range("H4").Value = '''no adverse impact on offsetting measures; compliance with other regulations and sara'
$ws.Range("I4").Value = '''authorization cannot be transferred without dfo notification'
$ws.Range("J4").Value = '''riverine, lake, shoreline'
$ws.Range("K4").Value = '''walleye'
$ws.Range("L4").Value = '''880 m2 (river habitat from old dam removal), 400 m2 (walleye spawning habitat), 895 m2 (shoreside/upperland), 720 m2 (lake to river habitat)'
$ws.Range("M4").Value = '''None'
$ws.Range("N4").Value = '''additional boulders to resist displacement if spawning areas wash out (contingency)'
$ws.Range("O4").Value = '''None'
$ws.Range("P4").Value = '''rock shoal (5-10m length, 10-40 m2 area)'

$ws.Range("A5").Value = '''OCR_18-HCAA-00064.json'
$ws.Range("B5").Value = '''None'
$ws.Range("C5").Value = '''None'
$ws.Range("D5").Value = '''None'
$ws.Range("E5").Value = '''Feb 07/2020'
$ws.Range("F5").Value = '''The Proponent must conduct monitoring and reporting on offsetting measures, including post-construction assessments and fish sampling, with reports due by June 14, 2022.'
$ws.Range("G5").Value = '''Structural stability and functionality of offsetting habitats must be maintained; contingency measures required if not met.'
$ws.Range("H5").Value = '''Prohibition on depositing deleterious substances in water frequented by fish.'
$ws.Range("I5").Value = '''Authorization cannot be transferred without prior notification to DFO.'
$ws.Range("J5").Value = '''Storm water management pond, tributary, drain, marsh'
$ws.Range("K5").Value = '''None explicitly listed in the provided sections'
$ws.Range("L5").Value = '''Storm water management pond (average depth 0.9m), unnamed tributary restoration, Hooper Drain channel, Central Drain channel, marsh habitat with berm'
$ws.Range("M5").Value = '''Riparian vegetation (trees and shrubs)'
$ws.Range("N5").Value = '''None mentioned'
$ws.Range("O5").Value = '''Habitat enhancement features (e.g., woody debris) in tributary restoration'
$ws.Range("P5").Value = '''Pool/deepwater habitat in SWM pond, riffles in tributary, channel morphology in drains'

$ws.Range("A6").Value = '''OCR_18-HCAA-00311.json'
$ws.Range("B6").Value = '''None'
$ws.Range("C6").Value = '''None'
$ws.Range("D6").Value = '''None'
$ws.Range("E6").Value = '''Oct 16, 2019'
$ws.Range("F6").Value = '''The Proponent must conduct annual monitoring reports on offsetting measures by Dec 31, 2020 and 2021, including photographic records and as-built surveys.'
$ws.Range("G6").Value = '''The Proponent must ensure no adverse impact on offsetting measures and obtain access permissions for DFO to monitor the measures.'
$ws.Range("H6").Value = '''The Proponent must implement erosion control measures, prevent deleterious substance deposits, and maintain a spill response plan.'
$ws.Range("I6").Value = '''The Proponent must adhere to the offsetting plan''s specifications, including vegetation survival rates and habitat utilization assessments.'
$ws.Range("J6").Value = '''Riverine'
$ws.Range("K6").Value = '''Rainbow Trout'
$ws.Range("L6").Value = '''Reach SN01: unspecified; Reach SN02: unspecified; unnamed tributary: unspecified'
$ws.Range("M6").Value = '''Riparian vegetation (e.g., trees, shrubs, and grass)'
$ws.Range("N6").Value = '''None explicitly mentioned'
$ws.Range("O6").Value = '''None explicitly mentioned'
$ws.Range("P6").Value = '''Pools, riffles, instream habitat features'

$ws.Range("A7").Value = '''OCR_18-HCAA-00192.json'
$ws.Range("B7").Value = '''53°36''40.96"N'
$ws.Range("C7").Value = '''108°44''38.01"W'
$ws.Range("D7").Value = '''53°36''40.96"N, 108°44''38.01"W; 12U 58311E E, 5940187 N'
$ws.Range("E7").Value = '''FEB 26 2020'
$ws.Range("F7").Value = '''Sediment and erosion control measures including installation of gravel/boulder bar during low flow periods.'
$ws.Range("G7").Value = '''Monitoring of offset structures post-construction for three years, including embeddedness surveys and sonar surveys.'
$ws.Range("H7").Value = '''Compliance with Species at Risk Act (SARA) restrictions, no harm to listed species or their habitats.'
$ws.Range("I7").Value = '''Submission of yearly reports (2022-2024) detailing monitoring results and adherence to offset criteria.'
$ws.Range("J7").Value = '''Riverine'
$ws.Range("K7").Value = '''None explicitly listed in section 4 or elsewhere'
$ws.Range("L7").Value = '''Gravel/boulder bar installation with area measurement required as per Offsetting Plan'
$ws.Range("M7").Value = '''Riparian vegetation (consultation with DFO on specifications)'
$ws.Range("N7").Value = '''Gravel/boulder bar installation as part of offset measures'
$ws.Range("O7").Value = '''None explicitly mentioned'
$ws.Range("P7").Value = '''Gravel/boulder bar (not extending beyond intake structure)'

$ws.Range("A8").Value = '''OCR_18-HCAA-00146.json'
$ws.Range("B8").Value = '''5672412N'
$ws.Range("C8").Value = '''11678490E'
$ws.Range("D8").Value = '''11678490E 5672412N'
$ws.Range("E8").Value = '''AUG 17 2018'
$ws.Range("F8").Value = '''Sedimentation and erosion control measures must be in place and maintained to avoid sediment release into the watercourse.'
$ws.Range("G8").Value = '''Total suspended sediment and turbidity monitoring must adhere to the plan in Appendix C of the 2018 Aquatic Effects Assessment.'
$ws.Range("H8").Value = '''All riprap must be clean and free of fine materials, not obtained from fish-frequented waters.'
$ws.Range("I8").Value = '''Contingency measures include upgrading erosion controls, increasing berm height, and halting work if turbidity thresholds are met.'
$ws.Range("J8").Value = '''Riverine'
$ws.Range("K8").Value = '''None explicitly listed in the document'
$ws.Range("L8").Value = '''7,800 m²'
$ws.Range("M8").Value = '''Well-vegetated area for dewatering discharge'
$ws.Range("N8").Value = '''Large boulders are to be stockpiled and replaced as part of natural structures preservation.'
$ws.Range("O8").Value = '''Woody debris is to be stockpiled and replaced as part of natural structures preservation.'
$ws.Range("P8").Value = '''Secondary channel re-graded to increase connectivity during low flow conditions.'

$ws.Range("A9").Value = '''OCR_18-HCAA-00145.json'
$ws.Range("B9").Value = '''43.79381'
$ws.Range("C9").Value = '''-80.386060'
$ws.Range("D9").Value = '''43.79381, -80.386060'
$ws.Range("E9").Value = '''JAN 10 2020'
$ws.Range("F9").Value = '''Sediment and erosion control measures must be in place and maintained to avoid sediment release.'
$ws.Range("G9").Value = '''Qualified environmental professional on-site to monitor instream and shoreline activities.'
$ws.Range("H9").Value = '''Fish rescue and relocation in isolated areas before work commences.'
$ws.Range("I9").Value = '''Temporary rock causeways must be installed no sooner than July 1 and removed by September 30 annually.'
$ws.Range("J9").Value = '''Riverine'
$ws.Range("K9").Value = '''trout'
$ws.Range("L9").Value = '''50 m2 boulder clusters, 100 m2 sweeper trees'
$ws.Range("M9").Value = '''Riparian vegetation (e.g. trees and shrubs and grass)'
$ws.Range("N9").Value = '''Boulder clusters for velocity refuge habitat'
$ws.Range("O9").Value = '''Anchored sweeper trees'
$ws.Range("P9").Value = '''Overwintering pool habitat with minimum 1.7m depth'

$ws.Range("A10").Value = '''OCR_18-HCAA-00253.json'
$ws.Range("B10").Value = '''50.894225'
$ws.Range("C10").Value = '''-114.009975'
$ws.Range("D10").Value = '''Longitude: -114.009975, Latitude: 50.894225'
$ws.Range("E10").Value = '''FEB 04 2018'
$ws.Range("F10").Value = '''Sediment and erosion control measures must be in place and maintained to avoid sediment release into water. Turbidity monitoring follows the 2018 plan. Contingency measures include upgrading erosion controls, adding armor, raising berms, and halting work if thresholds are met.'
$ws.Range("G10").Value = '''3,462 m² of habitat restoration upstream/downstream on the west bank per the offsetting plan. Monitoring reports due by Jan 15, 2023, with three years of post-construction monitoring ending Dec 31, 2025.'
$ws.Range("H10").Value = '''Prohibits deposit of deleterious substances, harming SARA-listed species, or damaging their habitats. Authorization cannot be transferred without DFO approval.'
$ws.Range("I10").Value = '''Proponent responsible for all design and safety aspects. Must comply with other regulatory agencies. Reports on mitigation measures post-work with photos and inspection reports.'
$ws.Range("J10").Value = '''Riverine'
$ws.Range("K10").Value = '''None explicitly listed in the document'
$ws.Range("L10").Value = '''3462 m²'
$ws.Range("M10").Value = '''None explicitly mentioned in offsetting measures'
$ws.Range("N10").Value = '''None explicitly mentioned'
$ws.Range("O10").Value = '''None explicitly mentioned'
$ws.Range("P10").Value = '''None explicitly mentioned in offsetting measures'

$ws.Range("A11").Value = '''OCR_18-HCAA-00160.json'
$ws.Range("B11").Value = '''71.889403°N'
$ws.Range("C11").Value = '''-80.887592°W'
$ws.Range("D11").Value = '''71.889403°N, -80.887592°W; Zone: 17 W, Easting: 503900 m E, Northings: 796600 m N'
$ws.Range("E11").Value = '''March 21, 2019'
$ws.Range("F11").Value = '''The work must be completed by the expiration date or DFO must be notified for extension.'
$ws.Range("G11").Value = '''Implement sediment and erosion control measures, including approved plans and monitoring turbidity levels.'
$ws.Range("H11").Value = '''Monitor and report on mitigation measures and submit reports by specified dates.'
$ws.Range("I11").Value = '''Offset habitat loss by placing course rock substrate to provide 2792 HEUs of habitat and implement contingency plans if needed.'
$ws.Range("J11").Value = '''Intertidal marine habitat, Subtidal marine habitat, Intertidal unnamed stream'
$ws.Range("K11").Value = '''None explicitly mentioned in the document'
$ws.Range("L11").Value = '''2792 HEUs of potential fish habitat'
$ws.Range("M11").Value = '''None explicitly mentioned'
$ws.Range("N11").Value = '''Course rock substrate placement as part of offset measures'
$ws.Range("O11").Value = '''None explicitly mentioned'
$ws.Range("P11").Value = '''None explicitly mentioned'
